# Diagramas de robustez y secuencia CU 16 y 17
# Update the "Estado" and "Esfuerzo (hrs)" columns for CU-16 (CRU gasto promocional, row 20)
# and CU-17 (CRU egreso, row 21) on the "Casos de Uso" sheet, moving both from
# "vacio" / 0 to "planificado" / 1, and move the active selection to E22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

# CU- 16 (row 20): Estado -> planificado, Esfuerzo (hrs) -> 1
$ws.Range("E20").Value = "planificado"
$ws.Range("F20").Value = 1

# CU- 17 (row 21): Estado -> planificado, Esfuerzo (hrs) -> 1
$ws.Range("E21").Value = "planificado"
$ws.Range("F21").Value = 1

# Update the active selection on the sheet to E22
$ws.Activate()
$ws.Range("E22").Select()
